$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: mark "Add exception handling to functions" as Complete
$ws.Range("C11").Value = "Complete"

# Copy formatting from the row above so the new row matches the table style
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)

# Add new row 12 for the new task
$ws.Range("A12").Value = Get-Date -Year 2025 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("B12").Value = "Implement price visualization functions"
$ws.Range("C12").Value = "Complete"
$ws.Range("D12").Value = "Adam Rodi"

# Move the selection to reflect where the user ended up
$ws.Range("B12").Select()
